# Applies the "Finalized - Mid Evaluation" text edits across slides 1-9.
$p = $ppt.ActivePresentation

# --- Slide 1: Title slide ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Deep Learning Explained"
$s1.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Artificial Intelligence • Machine Learning • Deep Learning"

# --- Slide 2: Lecture Agenda ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "3. Working of Neural Networks"

# --- Slide 3: Lecture Topics ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Artificial Intelligence"
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Deep Learning"
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).Runs(1).Text = "Applications"
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).Runs(1).Text = "Limitations"
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(7).Runs(1).Text = "Frameworks"

# --- Slide 4: Introduction to Deep Learning ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Deep learning powers features like Google Translate and phone gallery image grouping."
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "AI mimics human behavior, ML achieves AI through algorithms, and DL is ML inspired by the human brain."
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Deep learning utilizes artificial neural networks."

# --- Slide 5: Deep Learning vs. Machine Learning ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Example: Differentiating between tomatoes and cherries."
$s5.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Machine learning requires explicit feature definition (e.g., size, stem type)."
$s5.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Deep learning demands a much higher volume of data for training."

# --- Slide 6: How Neural Networks Work ---
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Working of Neural Networks"
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Neural networks can identify handwritten digits, represented as 28x28 pixel images."
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Each pixel is fed to a neuron in the input layer."
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Information is transferred through weighted channels and hidden layers to the output layer."
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Bias is added to the weighted sum of inputs, then applied to an activation function."
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).Runs(1).Text = "Activated neurons pass information, leading to the identification of the input digit."
$s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).Runs(1).Text = "Weights and bias are continuously adjusted to train the network."

# --- Slide 7: Applications of Deep Learning ---
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Self-driving cars: A growing reality with companies like Apple, Tesla, and Nissan."

# --- Slide 8: Limitations of Deep Learning ---
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Data: Requires a massive volume of data for effective training."
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Computational Power: Needs powerful and expensive Graphical Processing Units (GPUs)."
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Training Time: Can take hours or even months, increasing with data and network complexity."

# --- Slide 9: Frameworks and Future Outlook ---
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "The field of deep learning and AI is still in its early stages with vast future scope."
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Innovations like devices for the blind using deep learning and computer vision are emerging."
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Replicating the human mind may soon move beyond science fiction."
